$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap row 2 and row 3 values for columns D, J, K, L, M, P
$cols = @("D", "J", "K", "L", "M", "P")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")
    $tmp = $cell2.Value2
    $cell2.Value2 = $cell3.Value2
    $cell3.Value2 = $tmp
}
